# Realestate Update resale numbers 2024-01-22 22:08
# Appends a new data row (row 85) to the CityResaleNum sheet with the
# latest resale-number snapshot, mirroring the existing row layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 85

# Columns A (Date) and D (Week) hold text that *looks* like a date / a
# plain number ("2024-01-22", "03"). A bare .Value assignment would get
# auto-coerced into a date serial / numeric value, so we momentarily force
# a Text number format, assign the literal string, then clear the format
# again so the cell is left with no explicit style (matching the rest of
# the sheet) while the underlying value stays a plain text string.
$dateCell = $ws.Cells.Item($row, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2024-01-22"
$dateCell.ClearFormats()

$ws.Cells.Item($row, 2).Value = "22:08:45"
$ws.Cells.Item($row, 3).Value = "Monday"

$weekCell = $ws.Cells.Item($row, 4)
$weekCell.NumberFormat = "@"
$weekCell.Value = "03"
$weekCell.ClearFormats()

$ws.Cells.Item($row, 5).Value = 138498
$ws.Cells.Item($row, 6).Value = 141018
$ws.Cells.Item($row, 7).Value = 171361
$ws.Cells.Item($row, 8).Value = 148634
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 123218
$ws.Cells.Item($row, 11).Value = 223603
$ws.Cells.Item($row, 12).Value = 255950
$ws.Cells.Item($row, 13).Value = 185408
$ws.Cells.Item($row, 14).Value = 110292
$ws.Cells.Item($row, 15).Value = 41351
$ws.Cells.Item($row, 16).Value = 30892
$ws.Cells.Item($row, 17).Value = 73625
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42677
$ws.Cells.Item($row, 20).Value = -1
